$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgabenliste Projekt 1")

# Last task of the schedule ("Abgabe finales Video + Präsentation") is now
# complete -> set "% erledigt" to 100%. "Fortschritt" (F52) recomputes via
# its existing table formula.
$ws.Range("E52").Value = 1

# Author resized a few columns (E:H) while reviewing the finished sheet.
$ws.Columns.Item(5).ColumnWidth = 14.571428571428571
$ws.Columns.Item(6).ColumnWidth = 12.142857142857142
$ws.Columns.Item(7).ColumnWidth = 7
$ws.Columns.Item(8).ColumnWidth = 25.142857142857142

# Leave the cursor where the author ended up after saving.
$null = $ws.Range("E53").Select()
